# Add a WEEKEND_APPR_PROCESS_START column between ORGANIZATION_TYPE and
# FRAUD_RISK, derived from the existing WEEKDAY_APPR_PROCESS_START column
# (SATURDAY for weekend applications, N/A otherwise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column X (24) currently holds FRAUD_RISK; inserting here shifts it (and
# its header style) one column to the right, to Y, and leaves a fresh
# column X in its place.
$ws.Columns.Item(24).Insert()

$ws.Cells.Item(1, 24).Value = "WEEKEND_APPR_PROCESS_START"

$lastRow = $ws.Cells.Item(1, 21).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $weekday = $ws.Cells.Item($r, 21).Value()
    if ($weekday -eq "SATURDAY" -or $weekday -eq "SUNDAY") {
        $ws.Cells.Item($r, 24).Value = $weekday
    } else {
        $ws.Cells.Item($r, 24).Value = "N/A"
    }
}
